$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "2026-02-14T08:00:40.524701+00:00"
$ws.Range("E3").Value = "Please share ir"
$ws.Range("F3").Value = "yes"
